# Scrape refresh: 2026-01-12 01:41 JST
#
# A brand-new job listing ("HITOON" backend/payments work) showed up at
# the top of the freshly-scraped results, landing in row 4 and pushing
# the two rows that used to be there (UE5 / cordova) down by one. Every
# row's "取得日時" (fetched-at) timestamp is refreshed to the new scrape
# time, whether or not its other content changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-12 01:41:08"

# --- 1. Make room: push old rows 4-5 down to rows 5-6 -----------------
$ws.Rows.Item(4).Insert()

# --- 2. Refresh the retrieval timestamp on every data row -------------
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp

# --- 3. Write the newly scraped listing into row 4 ---------------------
$ws.Range("B4").Value = "フロント実装済み!音楽権利マーケットプレイス「HITOON」のバックエンド・決済機能実装"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5469298"
$ws.Range("G4").Value = 18

# --- 4. Rebuild the hyperlinks on column F ------------------------------
# Inserting the row shifted the cell *values* down automatically, but the
# worksheet's Hyperlinks collection keeps its old anchors (F2..F5), so it
# no longer lines up with the data. Clear it out and re-add one hyperlink
# per URL cell, now that the rows are in their final place.
$existingUrls = @(
    "https://www.lancers.jp/work/detail/5469128",
    "https://www.lancers.jp/work/detail/5468866",
    "https://www.lancers.jp/work/detail/5469203",
    "https://www.lancers.jp/work/detail/5469169"
)
foreach ($u in $existingUrls) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Address() -eq $u) {
            $h.Delete()
        }
    }
}

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5469128") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5468866") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5469298") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5469203") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5469169") | Out-Null

# Re-adding a hyperlink re-stamps the cell with a freshly minted style, so
# pin column F back onto the workbook's single shared "Hyperlink" style.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
